$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must stay plain text even though it parses as a
# number (e.g. "554.22"), without leaving a custom number format behind.
function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# --- Updated price / volume(1h) figures (and two reordered coin rows) ---
$ws.Range('D2').Value = '63.467.93'
$ws.Range('E2').Value = '  -1.00%  '
$ws.Range('D3').Value = '2.688.28'
$ws.Range('E3').Value = '  -2.51%  '
$ws.Range('E4').Value = '  -0.02%  '
Set-TextValue 'D5' '554.22'
$ws.Range('E5').Value = '  -3.75%  '
Set-TextValue 'D6' '158.06'
Set-TextValue 'D7' '0.999'
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('E8').Value = '  -3.24%  '
$ws.Range('E9').Value = '  -4.09%  '
$ws.Range('E11').Value = '  -4.43%  '
Set-TextValue 'D12' '5.37'
$ws.Range('E12').Value = '  -8.56%  '
$ws.Range('D13').Value = '3.167.53'
$ws.Range('E13').Value = '  -2.47%  '
Set-TextValue 'D14' '26.41'
$ws.Range('E14').Value = '  -2.01%  '
$ws.Range('D15').Value = '63.320.94'
$ws.Range('E15').Value = '  -0.80%  '
$ws.Range('D17').Value = '2.692.11'
$ws.Range('E17').Value = '  -2.65%  '
Set-TextValue 'D18' '12.07'
$ws.Range('E18').Value = '  -0.42%  '
$ws.Range('E19').Value = '  -4.80%  '
Set-TextValue 'D20' '343.19'
$ws.Range('E20').Value = '  -4.63%  '
Set-TextValue 'D21' '6.33'
$ws.Range('E21').Value = '  -4.61%  '
$ws.Range('E22').Value = '  -0.37%  '
$ws.Range('E23').Value = '  -3.98%  '
Set-TextValue 'D24' '63.89'
$ws.Range('E24').Value = '  -1.84%  '
Set-TextValue 'D26' '1.00'
$ws.Range('E26').Value = '  +0.14%  '
$ws.Range('E27').Value = '  -4.55%  '
$ws.Range('E28').Value = '  -5.03%  '
$ws.Range('E29').Value = '  -0.43%  '
Set-TextValue 'D30' '1.32'
$ws.Range('E30').Value = '  -3.01%  '
Set-TextValue 'D31' '7.01'
$ws.Range('E31').Value = '  -4.87%  '
Set-TextValue 'D32' '165.50'
Set-TextValue 'D33' '0.998'
$ws.Range('E33').Value = '  -0.05%  '
Set-TextValue 'D34' '4.78'
$ws.Range('E34').Value = '  -3.21%  '
Set-TextValue 'D35' '19.55'
$ws.Range('E35').Value = '  -3.33%  '
$ws.Range('E36').Value = '  -3.66%  '
$ws.Range('E37').Value = '  -1.81%  '
Set-TextValue 'D38' '340.10'
$ws.Range('E38').Value = '  -3.04%  '
Set-TextValue 'D39' '0.950'
$ws.Range('E39').Value = '  -5.58%  '
Set-TextValue 'D40' '6.05'
$ws.Range('E40').Value = '  -4.58%  '
Set-TextValue 'D41' '38.19'
$ws.Range('E41').Value = '  -2.29%  '
Set-TextValue 'D42' '3.93'
$ws.Range('E42').Value = '  -5.99%  '
Set-TextValue 'D43' '20.79'
$ws.Range('E43').Value = '  -5.51%  '
$ws.Range('B44').Value = 'Mantle'
$ws.Range('C44').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue 'D44' '0.622'
$ws.Range('E44').Value = '  -1.07%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D45' '20.30'
$ws.Range('E45').Value = '  -5.92%  '
Set-TextValue 'D46' '0.0563'
$ws.Range('E46').Value = '  -4.05%  '
$ws.Range('E47').Value = '  -0.06%  '
Set-TextValue 'D48' '11.08'
$ws.Range('E48').Value = '  +0.42%  '
Set-TextValue 'D49' '129.87'
$ws.Range('E49').Value = '  -5.45%  '
Set-TextValue 'D50' '0.0972'
$ws.Range('E50').Value = '  -3.91%  '
$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D51').Value = '2.101.78'
$ws.Range('E51').Value = '  -1.27%  '
